# Updated cryptos list values (price / 1h volume change) to match target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.294.93"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "2.097.71"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.69%  "
$ws.Range("D5").Value = "'342.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "'0.5277"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.4389"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "'55.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.45%  "
$ws.Range("D10").Value = "'0.09350"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Value = "'24.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "'8.565"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.08%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.112.78"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'6.865"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "'101.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "'0.00001158"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "'21.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "'0.06722"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "'6.383"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.89%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").Value = "30.276.64"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D24").Value = "'12.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("D25").Value = "'2.321"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "'7.007"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.54%  "
$ws.Range("D27").Value = "'21.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "'162.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "'2.513"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "'133.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "'1.134"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "'1.679"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").Value = "'6.243"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").Value = "'3.913"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("D36").Value = "'10.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("D37").Value = "'0.02621"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D38").Value = "'0.06761"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").Value = "'1.347"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").Value = "'0.6961"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "'0.6782"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'14.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'2.346"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "'1.002"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("E47").Value = "  +8.80%  "
$ws.Range("D48").Value = "'3.641"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "'0.00000000350"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").Value = "'1.211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.07%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.214"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
